$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set text (shared-string) values in the exact order the author must have
# typed them, so the shared-strings table ends up in the same order as
# the target workbook.
$ws.Range("B21").Value = "Multiplexer Select 1"
$ws.Range("B22").Value = "Multiplexer Select 2"
$ws.Range("B23").Value = "Multiplexer Select 3"

$ws.Range("H13").Value = "4051 Pinout"

$ws.Range("G14").Value = "CH 4 I/O"
$ws.Range("G15").Value = "CH 6 I/O"
$ws.Range("G16").Value = "COM O/I"
$ws.Range("G17").Value = "CH 7 I/O"
$ws.Range("G18").Value = "CH 5 I/O"

$ws.Range("G20").Value = "VEE"
$ws.Range("G21").Value = "VSS"

$ws.Range("J19").Value = "A"
$ws.Range("J20").Value = "B"
$ws.Range("J21").Value = "C"

$ws.Range("J15").Value = "CH 2 I/O"
$ws.Range("J16").Value = "CH 1 I/O"
$ws.Range("J17").Value = "CH 0 I/O"
$ws.Range("J18").Value = "CH 3 I/O"

$ws.Range("G19").Value = "INH"

$ws.Range("J14").Value = "VDD"

# Numeric pin-number cells (these don't go through the shared-strings table)
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 2
$ws.Range("H16").Value = 3
$ws.Range("H17").Value = 4
$ws.Range("H18").Value = 5
$ws.Range("H19").Value = 6
$ws.Range("H20").Value = 7
$ws.Range("H21").Value = 8

$ws.Range("I14").Value = 16
$ws.Range("I15").Value = 15
$ws.Range("I16").Value = 14
$ws.Range("I17").Value = 13
$ws.Range("I18").Value = 12
$ws.Range("I19").Value = 11
$ws.Range("I20").Value = 10
$ws.Range("I21").Value = 9

# Styling: shaded fill for the pin-number columns, with the pin column
# (H) additionally using a left-aligned integer number format.
$ws.Range("H14:H21").NumberFormat = "0"
$ws.Range("H14:H21").HorizontalAlignment = -4131
$ws.Range("H14:I21").Interior.ThemeColor = 1
$ws.Range("H14:I21").Interior.TintAndShade = 0.249977111117893

# Sheet view / selection, matching the final state in the workbook.
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("J25").Select()
